# Updated cryptos list values (Price / Volume(1h)) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.830.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "'2.481.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'556.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.19%  "
$ws.Range("D6").Value = "'148.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.21%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("D9").Value = "'2.480.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.93%  "
$ws.Range("E10").Value = "  -8.06%  "
$ws.Range("E11").Value = "  -5.07%  "
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("E13").Value = "  -6.42%  "
$ws.Range("D14").Value = "'26.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.15%  "
$ws.Range("D15").Value = "'2.930.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.88%  "
$ws.Range("E16").Value = "  -8.02%  "
$ws.Range("D17").Value = "'61.738.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.23%  "
$ws.Range("D18").Value = "'2.487.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.76%  "
$ws.Range("D19").Value = "'11.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.44%  "
$ws.Range("D20").Value = "'7.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.60%  "
$ws.Range("D21").Value = "'4.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.11%  "
$ws.Range("D22").Value = "'323.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.32%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "'64.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.23%  "
$ws.Range("E26").Value = "  -9.31%  "
$ws.Range("D27").Value = "'572.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "'2.609.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.90%  "
$ws.Range("E29").Value = "  -7.69%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -9.48%  "
$ws.Range("D32").Value = "'7.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.51%  "
$ws.Range("E33").Value = "  -6.15%  "
$ws.Range("E34").Value = "  -5.94%  "
$ws.Range("E35").Value = "  -6.42%  "
$ws.Range("E36").Value = "  -8.66%  "
$ws.Range("E37").Value = "  -8.98%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").Value = "'18.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("D41").Value = "'1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.77%  "
$ws.Range("D42").Value = "'144.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'2.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").Value = "'40.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("D46").Value = "'149.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.46%  "
$ws.Range("D47").Value = "'22.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.62%  "
$ws.Range("E48").Value = "  -6.32%  "
$ws.Range("D49").Value = "'0.0545"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.48%  "
$ws.Range("D50").Value = "'0.601"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.33%  "
$ws.Range("D51").Value = "'0.0945"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.60%  "
